$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(302, 44376, 0, 0, 0),
    @(303, 44377, 0, 0, 0),
    @(304, 44378, 0, 0, 0),
    @(305, 44379, 0, 0, 0),
    @(306, 44380, 0, 0, 0),
    @(307, 44381, 0, 0, 0),
    @(308, 44382, 0, 0, 0),
    @(309, 44383, 0, 0, 0),
    @(310, 44384, 0, 0, 0),
    @(311, 44385, 0, 0, 0),
    @(312, 44386, 0, 0, 0),
    @(313, 44387, 0, 0, 0),
    @(314, 44388, 0, 0, 0),
    @(315, 44389, 0, 0, 0),
    @(316, 44390, 1, 1, 23.82654276864427),
    @(317, 44391, 1, 2, 47.65308553728854),
    @(318, 44392, 1, 3, 71.47962830593281),
    @(319, 44393, 1, 4, 95.30617107457708),
    @(320, 44394, 0, 4, 95.30617107457708),
    @(321, 44395, 0, 4, 95.30617107457708),
    @(322, 44396, 0, 4, 95.30617107457708),
    @(323, 44397, 0, 3, 71.47962830593281),
    @(324, 44398, 0, 2, 47.65308553728854),
    @(325, 44399, 0, 1, 23.82654276864427),
    @(326, 44400, 1, 1, 23.82654276864427),
    @(327, 44401, 2, 3, 71.47962830593281),
    @(328, 44402, 1, 4, 95.30617107457708)
)

foreach ($item in $newData) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
}

# Copy the date-column number format/style down from the last existing row (A301)
# onto the newly added date cells (A302:A328), matching the style used throughout column A.
$ws.Range("A301").Copy()
$ws.Range("A302:A328").PasteSpecial(-4122)
$excel.CutCopyMode = 0
